$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 16.2726235
$ws.Range("H2").Value = 32.545247
$ws.Range("I2").Value = 0.1150896855633734
$ws.Range("J2").Value = 0.08350961674642161
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 67.073376
$ws.Range("N2").Value = 134.146752
$ws.Range("O2").Value = 0.3886152607163267
$ws.Range("P2").Value = 0.3201976973989085
$ws.Range("Q2").Value = 1091.459794521936
$ws.Range("R2").Value = 4365.839178087745
$ws.Range("S2").Value = 0.04472560816097042
$ws.Range("T2").Value = 0.02673958699286953
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 16.2726235
$ws.Range("H3").Value = 32.545247
$ws.Range("I3").Value = 0.1150896855633734
$ws.Range("J3").Value = 0.08350961674642161
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.219010666666667
$ws.Range("N3").Value = 12.657032
$ws.Range("O3").Value = 0.02444445215030819
$ws.Range("P3").Value = 0.03021133528677833
$ws.Range("Q3").Value = 68.65437212115067
$ws.Range("R3").Value = 411.9262327269041
$ws.Range("S3").Value = 0.002813304311747897
$ws.Range("T3").Value = 0.002522937031196502
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 16.2726235
$ws.Range("H4").Value = 32.545247
$ws.Range("I4").Value = 0.1150896855633734
$ws.Range("J4").Value = 0.08350961674642161
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.54988666666667
$ws.Range("N4").Value = 37.64966
$ws.Range("O4").Value = 0.07271256897710082
$ws.Range("P4").Value = 0.0898667635266472
$ws.Range("Q4").Value = 204.2195806943367
$ws.Range("R4").Value = 1225.31748416602
$ws.Range("S4").Value = 0.008368466700079635
$ws.Range("T4").Value = 0.007504738980351608
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 16.2726235
$ws.Range("H5").Value = 32.545247
$ws.Range("I5").Value = 0.1150896855633734
$ws.Range("J5").Value = 0.08350961674642161
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 12.18648933333333
$ws.Range("N5").Value = 36.559468
$ws.Range("O5").Value = 0.07060708751994335
$ws.Range("P5").Value = 0.08726456136432641
$ws.Range("Q5").Value = 198.3061527080993
$ws.Range("R5").Value = 1189.836916248596
$ws.Range("S5").Value = 0.008126147501215868
$ws.Range("T5").Value = 0.007287430075079489
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 16.2726235
$ws.Range("H6").Value = 32.545247
$ws.Range("I6").Value = 0.1150896855633734
$ws.Range("J6").Value = 0.08350961674642161
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.80271133333334
$ws.Range("N6").Value = 134.408134
$ws.Range("O6").Value = 0.2595816460111038
$ws.Range("P6").Value = 0.3208215955797718
$ws.Range("Q6").Value = 729.0576533065165
$ws.Range("R6").Value = 4374.345919839099
$ws.Range("S6").Value = 0.02987517001744084
$ws.Range("T6").Value = 0.02679168849084222
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 16.2726235
$ws.Range("H7").Value = 32.545247
$ws.Range("I7").Value = 0.1150896855633734
$ws.Range("J7").Value = 0.08350961674642161
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 31.7643625
$ws.Range("N7").Value = 63.528725
$ws.Range("O7").Value = 0.1840389846252172
$ws.Range("P7").Value = 0.1516380468435678
$ws.Range("Q7").Value = 516.8895116800188
$ws.Range("R7").Value = 2067.558046720075
$ws.Range("S7").Value = 0.02118098887191876
$ws.Range("T7").Value = 0.01266323517608227
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.93153
$ws.Range("H8").Value = 71.79459
$ws.Range("I8").Value = 0.169258034068719
$ws.Range("J8").Value = 0.1842216375056694
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 67.073376
$ws.Range("N8").Value = 134.146752
$ws.Range("O8").Value = 0.3886152607163267
$ws.Range("P8").Value = 0.3201976973989085
$ws.Range("Q8").Value = 1605.16850994528
$ws.Range("R8").Value = 9631.011059671679
$ws.Range("S8").Value = 0.06577625503794815
$ws.Range("T8").Value = 0.05898734414037173
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.93153
$ws.Range("H9").Value = 71.79459
$ws.Range("I9").Value = 0.169258034068719
$ws.Range("J9").Value = 0.1842216375056694
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.219010666666667
$ws.Range("N9").Value = 12.657032
$ws.Range("O9").Value = 0.02444445215030819
$ws.Range("P9").Value = 0.03021133528677833
$ws.Range("Q9").Value = 100.9673803396533
$ws.Range("R9").Value = 908.7064230568801
$ws.Range("S9").Value = 0.004137419914848036
$ws.Range("T9").Value = 0.005565581657763115
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.93153
$ws.Range("H10").Value = 71.79459
$ws.Range("I10").Value = 0.169258034068719
$ws.Range("J10").Value = 0.1842216375056694
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.54988666666667
$ws.Range("N10").Value = 37.64966
$ws.Range("O10").Value = 0.07271256897710082
$ws.Range("P10").Value = 0.0898667635266472
$ws.Range("Q10").Value = 300.3379892599333
$ws.Range("R10").Value = 2703.0419033394
$ws.Range("S10").Value = 0.01230718647715021
$ws.Range("T10").Value = 0.01655540233421371
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 23.93153
$ws.Range("H11").Value = 71.79459
$ws.Range("I11").Value = 0.169258034068719
$ws.Range("J11").Value = 0.1842216375056694
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 12.18648933333333
$ws.Range("N11").Value = 36.559468
$ws.Range("O11").Value = 0.07060708751994335
$ws.Range("P11").Value = 0.08726456136432641
$ws.Range("Q11").Value = 291.6413350753466
$ws.Range("R11").Value = 2624.77201567812
$ws.Range("S11").Value = 0.0119508168249436
$ws.Range("T11").Value = 0.01607602039075018
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 23.93153
$ws.Range("H12").Value = 71.79459
$ws.Range("I12").Value = 0.169258034068719
$ws.Range("J12").Value = 0.1842216375056694
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.80271133333334
$ws.Range("N12").Value = 134.408134
$ws.Range("O12").Value = 0.2595816460111038
$ws.Range("P12").Value = 0.3208215955797718
$ws.Range("Q12").Value = 1072.197430355007
$ws.Range("R12").Value = 9649.776873195062
$ws.Range("S12").Value = 0.04393627908416156
$ws.Range("T12").Value = 0.05910227968488718
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 23.93153
$ws.Range("H13").Value = 71.79459
$ws.Range("I13").Value = 0.169258034068719
$ws.Range("J13").Value = 0.1842216375056694
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 31.7643625
$ws.Range("N13").Value = 63.528725
$ws.Range("O13").Value = 0.1840389846252172
$ws.Range("P13").Value = 0.1516380468435678
$ws.Range("Q13").Value = 760.1697940996249
$ws.Range("R13").Value = 4561.01876459775
$ws.Range("S13").Value = 0.03115007672966747
$ws.Range("T13").Value = 0.02793500929768346
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 35.08101266666667
$ws.Range("H14").Value = 105.243038
$ws.Range("I14").Value = 0.2481138162541146
$ws.Range("J14").Value = 0.2700488267490822
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 67.073376
$ws.Range("N14").Value = 134.146752
$ws.Range("O14").Value = 0.3886152607163267
$ws.Range("P14").Value = 0.3201976973989085
$ws.Range("Q14").Value = 2353.001953052096
$ws.Range("R14").Value = 14118.01171831258
$ws.Range("S14").Value = 0.09642081539091553
$ws.Range("T14").Value = 0.08646901251033288
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 35.08101266666667
$ws.Range("H15").Value = 105.243038
$ws.Range("I15").Value = 0.2481138162541146
$ws.Range("J15").Value = 0.2700488267490822
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.219010666666667
$ws.Range("N15").Value = 12.657032
$ws.Range("O15").Value = 0.02444445215030819
$ws.Range("P15").Value = 0.03021133528677833
$ws.Range("Q15").Value = 148.0071666381351
$ws.Range("R15").Value = 1332.064499743216
$ws.Range("S15").Value = 0.006065006309254064
$ws.Range("T15").Value = 0.008158535648717635
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 35.08101266666667
$ws.Range("H16").Value = 105.243038
$ws.Range("I16").Value = 0.2481138162541146
$ws.Range("J16").Value = 0.2700488267490822
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.54988666666667
$ws.Range("N16").Value = 37.64966
$ws.Range("O16").Value = 0.07271256897710082
$ws.Range("P16").Value = 0.0898667635266472
$ws.Range("Q16").Value = 440.2627331185645
$ws.Range("R16").Value = 3962.36459806708
$ws.Range("S16").Value = 0.01804099297854903
$ws.Range("T16").Value = 0.02426841405410829
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 35.08101266666667
$ws.Range("H17").Value = 105.243038
$ws.Range("I17").Value = 0.2481138162541146
$ws.Range("J17").Value = 0.2700488267490822
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 12.18648933333333
$ws.Range("N17").Value = 36.559468
$ws.Range("O17").Value = 0.07060708751994335
$ws.Range("P17").Value = 0.08726456136432641
$ws.Range("Q17").Value = 427.514386664865
$ws.Range("R17").Value = 3847.629479983784
$ws.Range("S17").Value = 0.01751859393916141
$ws.Range("T17").Value = 0.02356569241320964
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 35.08101266666667
$ws.Range("H18").Value = 105.243038
$ws.Range("I18").Value = 0.2481138162541146
$ws.Range("J18").Value = 0.2700488267490822
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 44.80271133333334
$ws.Range("N18").Value = 134.408134
$ws.Range("O18").Value = 0.2595816460111038
$ws.Range("P18").Value = 0.3208215955797718
$ws.Range("Q18").Value = 1571.724483785677
$ws.Range("R18").Value = 14145.5203540711
$ws.Range("S18").Value = 0.06440579282133962
$ws.Range("T18").Value = 0.08663749548208591
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = 35.08101266666667
$ws.Range("H19").Value = 105.243038
$ws.Range("I19").Value = 0.2481138162541146
$ws.Range("J19").Value = 0.2700488267490822
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 31.7643625
$ws.Range("N19").Value = 63.528725
$ws.Range("O19").Value = 0.1840389846252172
$ws.Range("P19").Value = 0.1516380468435678
$ws.Range("Q19").Value = 1114.326003211092
$ws.Range("R19").Value = 6685.956019266551
$ws.Range("S19").Value = 0.04566261481489496
$ws.Range("T19").Value = 0.04094967664062785
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 37.788217
$ws.Range("H20").Value = 113.364651
$ws.Range("I20").Value = 0.2672607777430924
$ws.Range("J20").Value = 0.290888514614802
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 67.073376
$ws.Range("N20").Value = 134.146752
$ws.Range("O20").Value = 0.3886152607163267
$ws.Range("P20").Value = 0.3201976973989085
$ws.Range("Q20").Value = 2534.583287210592
$ws.Range("R20").Value = 15207.49972326355
$ws.Range("S20").Value = 0.1038616168218801
$ws.Range("T20").Value = 0.09314183257944834
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 37.788217
$ws.Range("H21").Value = 113.364651
$ws.Range("I21").Value = 0.2672607777430924
$ws.Range("J21").Value = 0.290888514614802
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 4.219010666666667
$ws.Range("N21").Value = 12.657032
$ws.Range("O21").Value = 0.02444445215030819
$ws.Range("P21").Value = 0.03021133528677833
$ws.Range("Q21").Value = 159.4288905973147
$ws.Range("R21").Value = 1434.860015375832
$ws.Range("S21").Value = 0.006533043293195175
$ws.Range("T21").Value = 0.008788130446100701
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 37.788217
$ws.Range("H22").Value = 113.364651
$ws.Range("I22").Value = 0.2672607777430924
$ws.Range("J22").Value = 0.290888514614802
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 12.54988666666667
$ws.Range("N22").Value = 37.64966
$ws.Range("O22").Value = 0.07271256897710082
$ws.Range("P22").Value = 0.0898667635266472
$ws.Range("Q22").Value = 474.2378406854067
$ws.Range("R22").Value = 4268.14056616866
$ws.Range("S22").Value = 0.01943321773651822
$ws.Range("T22").Value = 0.02614120935550607
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 37.788217
$ws.Range("H23").Value = 113.364651
$ws.Range("I23").Value = 0.2672607777430924
$ws.Range("J23").Value = 0.290888514614802
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 12.18648933333333
$ws.Range("N23").Value = 36.559468
$ws.Range("O23").Value = 0.07060708751994335
$ws.Range("P23").Value = 0.08726456136432641
$ws.Range("Q23").Value = 460.5057033961853
$ws.Range("R23").Value = 4144.551330565668
$ws.Range("S23").Value = 0.01887050512475465
$ws.Range("T23").Value = 0.02538425863378115
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 37.788217
$ws.Range("H24").Value = 113.364651
$ws.Range("I24").Value = 0.2672607777430924
$ws.Range("J24").Value = 0.290888514614802
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 44.80271133333334
$ws.Range("N24").Value = 134.408134
$ws.Range("O24").Value = 0.2595816460111038
$ws.Range("P24").Value = 0.3208215955797718
$ws.Range("Q24").Value = 1693.01457805236
$ws.Range("R24").Value = 15237.13120247124
$ws.Range("S24").Value = 0.06937599260075969
$ws.Range("T24").Value = 0.09332331739455056
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 37.788217
$ws.Range("H25").Value = 113.364651
$ws.Range("I25").Value = 0.2672607777430924
$ws.Range("J25").Value = 0.290888514614802
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 31.7643625
$ws.Range("N25").Value = 63.528725
$ws.Range("O25").Value = 0.1840389846252172
$ws.Range("P25").Value = 0.1516380468435678
$ws.Range("Q25").Value = 1200.318623016663
$ws.Range("R25").Value = 7201.911738099976
$ws.Range("S25").Value = 0.04918640216598457
$ws.Range("T25").Value = 0.04410976620541521
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 10.13617966666667
$ws.Range("H26").Value = 30.408539
$ws.Range("I26").Value = 0.07168909983387287
$ws.Range("J26").Value = 0.07802692164876225
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 67.073376
$ws.Range("N26").Value = 134.146752
$ws.Range("O26").Value = 0.3886152607163267
$ws.Range("P26").Value = 0.3201976973989085
$ws.Range("Q26").Value = 679.867789985888
$ws.Range("R26").Value = 4079.206739915328
$ws.Range("S26").Value = 0.02785947822245928
$ws.Range("T26").Value = 0.02498404064705871
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 10.13617966666667
$ws.Range("H27").Value = 30.408539
$ws.Range("I27").Value = 0.07168909983387287
$ws.Range("J27").Value = 0.07802692164876225
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 4.219010666666667
$ws.Range("N27").Value = 12.657032
$ws.Range("O27").Value = 0.02444445215030819
$ws.Range("P27").Value = 0.03021133528677833
$ws.Range("Q27").Value = 42.76465013291644
$ws.Range("R27").Value = 384.881851196248
$ws.Range("S27").Value = 0.001752400770587773
$ws.Range("T27").Value = 0.002357297491325939
$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 10.13617966666667
$ws.Range("H28").Value = 30.408539
$ws.Range("I28").Value = 0.07168909983387287
$ws.Range("J28").Value = 0.07802692164876225
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 12.54988666666667
$ws.Range("N28").Value = 37.64966
$ws.Range("O28").Value = 0.07271256897710082
$ws.Range("P28").Value = 0.0898667635266472
$ws.Range("Q28").Value = 127.2079060496378
$ws.Range("R28").Value = 1144.87115444674
$ws.Range("S28").Value = 0.005212698616576748
$ws.Range("T28").Value = 0.007012026916521546
$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 10.13617966666667
$ws.Range("H29").Value = 30.408539
$ws.Range("I29").Value = 0.07168909983387287
$ws.Range("J29").Value = 0.07802692164876225
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 12.18648933333333
$ws.Range("N29").Value = 36.559468
$ws.Range("O29").Value = 0.07060708751994335
$ws.Range("P29").Value = 0.08726456136432641
$ws.Range("Q29").Value = 123.5244453885835
$ws.Range("R29").Value = 1111.720008497252
$ws.Range("S29").Value = 0.005061758546196218
$ws.Range("T29").Value = 0.006808985092287902
$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 10.13617966666667
$ws.Range("H30").Value = 30.408539
$ws.Range("I30").Value = 0.07168909983387287
$ws.Range("J30").Value = 0.07802692164876225
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 44.80271133333334
$ws.Range("N30").Value = 134.408134
$ws.Range("O30").Value = 0.2595816460111038
$ws.Range("P30").Value = 0.3208215955797718
$ws.Range("Q30").Value = 454.1283316284697
$ws.Range("R30").Value = 4087.154984656227
$ws.Range("S30").Value = 0.01860917453593107
$ws.Range("T30").Value = 0.02503272150153374
$ws.Range("E31").Value = 2
$ws.Range("G31").Value = 10.13617966666667
$ws.Range("H31").Value = 30.408539
$ws.Range("I31").Value = 0.07168909983387287
$ws.Range("J31").Value = 0.07802692164876225
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 31.7643625
$ws.Range("N31").Value = 63.528725
$ws.Range("O31").Value = 0.1840389846252172
$ws.Range("P31").Value = 0.1516380468435678
$ws.Range("Q31").Value = 321.9692852971292
$ws.Range("R31").Value = 1931.815711782775
$ws.Range("S31").Value = 0.01319358914212179
$ws.Range("T31").Value = 0.0118318500000344
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 18.181244
$ws.Range("H32").Value = 36.362488
$ws.Range("I32").Value = 0.1285885865368279
$ws.Range("J32").Value = 0.09330448273526253
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 67.073376
$ws.Range("N32").Value = 134.146752
$ws.Range("O32").Value = 0.3886152607163267
$ws.Range("P32").Value = 0.3201976973989085
$ws.Range("Q32").Value = 1219.477414959744
$ws.Range("R32").Value = 4877.909659838976
$ws.Range("S32").Value = 0.04997148708215331
$ws.Range("T32").Value = 0.02987588052882727
$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 18.181244
$ws.Range("H33").Value = 36.362488
$ws.Range("I33").Value = 0.1285885865368279
$ws.Range("J33").Value = 0.09330448273526253
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 4.219010666666667
$ws.Range("N33").Value = 12.657032
$ws.Range("O33").Value = 0.02444445215030819
$ws.Range("P33").Value = 0.03021133528677833
$ws.Range("Q33").Value = 76.70686236926933
$ws.Range("R33").Value = 460.241174215616
$ws.Range("S33").Value = 0.003143277550675254
$ws.Range("T33").Value = 0.002818853011674436
$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 18.181244
$ws.Range("H34").Value = 36.362488
$ws.Range("I34").Value = 0.1285885865368279
$ws.Range("J34").Value = 0.09330448273526253
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 12.54988666666667
$ws.Range("N34").Value = 37.64966
$ws.Range("O34").Value = 0.07271256897710082
$ws.Range("P34").Value = 0.0898667635266472
$ws.Range("Q34").Value = 228.1725516590133
$ws.Range("R34").Value = 1369.03530995408
$ws.Range("S34").Value = 0.009350006468226996
$ws.Range("T34").Value = 0.008384971885945974
$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 18.181244
$ws.Range("H35").Value = 36.362488
$ws.Range("I35").Value = 0.1285885865368279
$ws.Range("J35").Value = 0.09330448273526253
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 12.18648933333333
$ws.Range("N35").Value = 36.559468
$ws.Range("O35").Value = 0.07060708751994335
$ws.Range("P35").Value = 0.08726456136432641
$ws.Range("Q35").Value = 221.5655360727307
$ws.Range("R35").Value = 1329.393216436384
$ws.Range("S35").Value = 0.009079265583671616
$ws.Range("T35").Value = 0.00814217475921805
$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 18.181244
$ws.Range("H36").Value = 36.362488
$ws.Range("I36").Value = 0.1285885865368279
$ws.Range("J36").Value = 0.09330448273526253
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 44.80271133333334
$ws.Range("N36").Value = 134.408134
$ws.Range("O36").Value = 0.2595816460111038
$ws.Range("P36").Value = 0.3208215955797718
$ws.Range("Q36").Value = 814.5690266128988
$ws.Range("R36").Value = 4887.414159677392
$ws.Range("S36").Value = 0.03337923695147104
$ws.Range("T36").Value = 0.0299340930258722
$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 18.181244
$ws.Range("H37").Value = 36.362488
$ws.Range("I37").Value = 0.1285885865368279
$ws.Range("J37").Value = 0.09330448273526253
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 31.7643625
$ws.Range("N37").Value = 63.528725
$ws.Range("O37").Value = 0.1840389846252172
$ws.Range("P37").Value = 0.1516380468435678
$ws.Range("Q37").Value = 577.51562511695
$ws.Range("R37").Value = 2310.0625004678
$ws.Range("S37").Value = 0.02366531290062968
$ws.Range("T37").Value = 0.0141485095237246
